$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Corrected Avg_Time_ms measurements for the 5,000- and 10,000-row runs
# (re-run after fixing the sort routine).
$ws.Range("D2").Value = 0.94876558
$ws.Range("D3").Value = 2.0802711
